$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.182.29'
$ws.Range('E2').Value = '  -0.47%  '

$ws.Range('D3').Value = '2.641.25'
$ws.Range('E3').Value = '  -0.29%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').Value = '''594.41'
$ws.Range('E5').Value = '  -0.84%  '

$ws.Range('D6').Value = '''158.93'
$ws.Range('E6').Value = '  +2.73%  '

$ws.Range('E7').Value = '  +0.00%  '

$ws.Range('D8').Value = '''0.543'
$ws.Range('E8').Value = '  -0.97%  '

$ws.Range('E9').Value = '  -2.67%  '

$ws.Range('E10').Value = '  -1.48%  '

$ws.Range('E11').Value = '  -0.32%  '

$ws.Range('E12').Value = '  -1.68%  '

$ws.Range('E13').Value = '  -1.78%  '

$ws.Range('D14').Value = '3.121.08'
$ws.Range('E14').Value = '  -0.30%  '

$ws.Range('D15').Value = '''0.0000187'
$ws.Range('E15').Value = '  -3.44%  '

$ws.Range('D16').Value = '68.023.11'
$ws.Range('E16').Value = '  -0.64%  '

$ws.Range('D17').Value = '2.611.29'
$ws.Range('E17').Value = '  -1.16%  '

$ws.Range('D18').Value = '''11.33'
$ws.Range('E18').Value = '  -1.48%  '

$ws.Range('D19').Value = '''359.63'
$ws.Range('E19').Value = '  -1.92%  '

$ws.Range('D20').Value = '''7.32'
$ws.Range('E20').Value = '  -2.65%  '

$ws.Range('D21').Value = '''4.38'
$ws.Range('E21').Value = '  -0.69%  '

$ws.Range('E22').Value = '  -3.65%  '

$ws.Range('E23').Value = '  -0.97%  '

$ws.Range('D24').Value = '''74.77'
$ws.Range('E24').Value = '  +1.32%  '

$ws.Range('E25').Value = '  +0.09%  '

$ws.Range('D26').Value = '''9.75'
$ws.Range('E26').Value = '  -1.80%  '

$ws.Range('E27').Value = '  -0.27%  '

$ws.Range('E28').Value = '  -4.40%  '

$ws.Range('E29').Value = '  +0.05%  '

$ws.Range('D30').Value = '''560.33'
$ws.Range('E30').Value = '  -3.57%  '

$ws.Range('D31').Value = '''7.98'
$ws.Range('E31').Value = '  -3.05%  '

$ws.Range('D32').Value = '''1.38'
$ws.Range('E32').Value = '  -4.34%  '

$ws.Range('E33').Value = '  -1.54%  '

$ws.Range('E34').Value = '  +0.01%  '

$ws.Range('E35').Value = '  -3.76%  '

$ws.Range('D36').Value = '''1.55'
$ws.Range('E36').Value = '  -4.29%  '

$ws.Range('B37').Value = 'EthereumClassic'
$ws.Range('C37').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D37').Value = '''19.67'
$ws.Range('E37').Value = '  +0.90%  '

$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D38').Value = '''158.79'
$ws.Range('E38').Value = '  -0.49%  '

$ws.Range('E39').Value = '  -1.72%  '

$ws.Range('D40').Value = '''1.86'
$ws.Range('E40').Value = '  -2.66%  '

$ws.Range('D41').Value = '''5.29'
$ws.Range('E41').Value = '  -3.11%  '

$ws.Range('D43').Value = '''2.60'
$ws.Range('E43').Value = '  -3.97%  '

$ws.Range('D44').Value = '0.0₆0322'
$ws.Range('E44').Value = '  -4.33%  '

$ws.Range('E45').Value = '  +0.02%  '

$ws.Range('D46').Value = '''156.81'
$ws.Range('E46').Value = '  -0.70%  '

$ws.Range('D47').Value = '''3.73'
$ws.Range('E47').Value = '  -2.03%  '

$ws.Range('D48').Value = '''21.84'
$ws.Range('E48').Value = '  -0.95%  '

$ws.Range('E49').Value = '  -2.79%  '

$ws.Range('E50').Value = '  -2.25%  '

$ws.Range('D51').Value = '''0.611'
$ws.Range('E51').Value = '  -1.09%  '
